$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 9273.429
$ws.Range("I5").Value = 8433.799999999999
$ws.Range("J5").Value = 11372.5
$ws.Range("K5").Value = 8433.799999999999
$ws.Range("L5").Value = 11372.5
$ws.Range("M5").Value = -8318.799999999999
$ws.Range("N5").Value = -11602.5
$ws.Range("H17").Value = 3963
$ws.Range("J17").Value = 3963
$ws.Range("L17").Value = 11889
$ws.Range("N17").Value = -12225
$ws.Range("H32").Value = 2248.8572
$ws.Range("I32").Value = 1165.6666
$ws.Range("J32").Value = 3061.25
$ws.Range("K32").Value = 1165.6666
$ws.Range("L32").Value = 3061.25
$ws.Range("M32").Value = -839.6666
$ws.Range("N32").Value = -3713.25
$ws.Range("H107").Value = 506.92856
$ws.Range("J107").Value = 758.1667
$ws.Range("L107").Value = 758.1667
$ws.Range("N107").Value = -4598.1667
$ws.Range("H138").Value = 2194.25
$ws.Range("J138").Value = 2647.791
$ws.Range("L138").Value = 7943.373000000001
$ws.Range("N138").Value = -18223.373
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
$ws.Range("H140").Value = 69873.11
$ws.Range("J140").Value = 69768.625
$ws.Range("L140").Value = 69768.625
$ws.Range("N140").Value = -80128.625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 730125.5600000001
$ws.Range("I2").Value = 796282.4399999999
$ws.Range("J2").Value = 2400
$ws.Range("K2").Value = 796282.4399999999
$ws.Range("L2").Value = 2400
$ws.Range("M2").Value = -796169.4399999999
$ws.Range("N2").Value = -2626
$ws.Range("H32").Value = 20969.174
$ws.Range("I32").Value = 22917.412
$ws.Range("K32").Value = 22917.412
$ws.Range("M32").Value = -22630.412
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16608
$ws.Range("H102").Value = 443197.97
$ws.Range("I102").Value = 527847.5600000001
$ws.Range("K102").Value = 527847.5600000001
$ws.Range("M102").Value = -526225.5600000001
$ws.Range("H116").Value = 730125.5600000001
$ws.Range("I116").Value = 796282.4399999999
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 796282.4399999999
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = -793988.4399999999
$ws.Range("N116").Value = -6988
$ws.Range("H132").Value = 11807.707
$ws.Range("I132").Value = 13334.404
$ws.Range("K132").Value = 40003.212
$ws.Range("M132").Value = -37473.212

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 730125.5600000001
$ws.Range("I3").Value = 796282.4399999999
$ws.Range("J3").Value = 2400
$ws.Range("K3").Value = 796282.4399999999
$ws.Range("L3").Value = 2400
$ws.Range("M3").Value = -796168.4399999999
$ws.Range("N3").Value = -2628
$ws.Range("H22").Value = 454.5
$ws.Range("I22").Value = 484.375
$ws.Range("K22").Value = 484.375
$ws.Range("M22").Value = -311.375
$ws.Range("H75").Value = 100236
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 100236
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H86").Value = 4009.4443
$ws.Range("I86").Value = 3809.4
$ws.Range("K86").Value = 3809.4
$ws.Range("M86").Value = -2686.4
$ws.Range("H88").Value = 25547
$ws.Range("J88").Value = 20896
$ws.Range("L88").Value = 20896
$ws.Range("N88").Value = -21708
$ws.Range("H89").Value = 4009.4443
$ws.Range("I89").Value = 3809.4
$ws.Range("K89").Value = 19047
$ws.Range("M89").Value = -13431
$ws.Range("H91").Value = 25547
$ws.Range("J91").Value = 20896
$ws.Range("L91").Value = 20896
$ws.Range("N91").Value = -23704
$ws.Range("H134").Value = 1147.4073
$ws.Range("I134").Value = 1114.6154
$ws.Range("K134").Value = 3343.8462
$ws.Range("M134").Value = -808.8462

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 251.84616
$ws.Range("I7").Value = 196.16667
$ws.Range("J7").Value = 299.57144
$ws.Range("K7").Value = 196.16667
$ws.Range("L7").Value = 299.57144
$ws.Range("M7").Value = -83.16667000000001
$ws.Range("N7").Value = -525.5714399999999
$ws.Range("H16").Value = 2820
$ws.Range("I16").Value = 2825
$ws.Range("J16").Value = 2800
$ws.Range("K16").Value = 2825
$ws.Range("L16").Value = 2800
$ws.Range("M16").Value = -2538
$ws.Range("N16").Value = -3374
$ws.Range("H31").Value = 6694.5806
$ws.Range("I31").Value = 3091.875
$ws.Range("K31").Value = 3091.875
$ws.Range("M31").Value = -2796.875
$ws.Range("H34").Value = 6694.5806
$ws.Range("I34").Value = 3091.875
$ws.Range("K34").Value = 3091.875
$ws.Range("M34").Value = -2889.875
$ws.Range("H58").Value = 576262.7
$ws.Range("I58").Value = 627137.2
$ws.Range("K58").Value = 627137.2
$ws.Range("M58").Value = -626934.2
$ws.Range("H105").Value = 5683193
$ws.Range("I105").Value = 22727272
$ws.Range("K105").Value = 22727272
$ws.Range("M105").Value = -22725525
$ws.Range("H107").Value = 1399079.4
$ws.Range("I107").Value = 1818551.9
$ws.Range("K107").Value = 1818551.9
$ws.Range("M107").Value = -1816631.9
$ws.Range("H113").Value = 2820
$ws.Range("I113").Value = 2825
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 2825
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -655
$ws.Range("N113").Value = -7140
$ws.Range("H136").Value = 576262.7
$ws.Range("I136").Value = 627137.2
$ws.Range("K136").Value = 1881411.6
$ws.Range("M136").Value = -1878861.6
$ws.Range("H141").Value = 73093.25
$ws.Range("I141").Value = 58696
$ws.Range("J141").Value = 77892.336
$ws.Range("K141").Value = 58696
$ws.Range("L141").Value = 77892.336
$ws.Range("M141").Value = -53516
$ws.Range("N141").Value = -88252.336

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 211.3077
$ws.Range("I14").Value = 211.3077
$ws.Range("K14").Value = 633.9231
$ws.Range("M14").Value = -460.9231
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("N41").Value = 0
$ws.Range("H97").Value = 584
$ws.Range("I97").Value = 376
$ws.Range("K97").Value = 1128
$ws.Range("M97").Value = -632

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9524169
$ws.Range("I107").Value = 11905112
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 11905112
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = -11903192
$ws.Range("N107").Value = -4240
$ws.Range("H113").Value = 2995
$ws.Range("I113").Value = 1990
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1990
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 180
$ws.Range("N113").Value = -8340
$ws.Range("H140").Value = 29999.5
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 89998
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 89998
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -100358

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 16930
$ws.Range("I32").Value = 16930
$ws.Range("K32").Value = 16930
$ws.Range("M32").Value = -16613
$ws.Range("H46").Value = 5765.269
$ws.Range("I46").Value = 2900
$ws.Range("K46").Value = 2900
$ws.Range("M46").Value = -2712
$ws.Range("H132").Value = 3560.3242
$ws.Range("I132").Value = 3141.1
$ws.Range("K132").Value = 9423.299999999999
$ws.Range("M132").Value = -6893.299999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1079879.5
$ws.Range("I100").Value = 1438450.4
$ws.Range("K100").Value = 2876900.8
$ws.Range("M100").Value = -2876359.8
$ws.Range("H107").Value = 3805.2273
$ws.Range("I107").Value = 4569.6924
$ws.Range("J107").Value = 2701
$ws.Range("K107").Value = 13709.0772
$ws.Range("L107").Value = 8103
$ws.Range("M107").Value = -11789.0772
$ws.Range("N107").Value = -11943
